$d = $word.ActiveDocument

# The last paragraph in the document body (currently empty) gets filled
# with the new sentence, split across many runs (each run holding one
# of the original "segments" so the resulting OOXML <w:r> boundaries
# match the authored edit).
$n = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($n)
$start0 = $p.Range.Start

$parts = @(
  " ",
  "Target System",
  "에서 무언가 스크립트를 다운받아 실행하고 싶다면 ",
  "tmp directory",
  "에 설치하는 편이 좋다",
  ". Tmp directory",
  "의 경우 어떤 ",
  "user",
  "던 다운로드와 파일을 실행할 수 있다",
  ".(",
  "아마도 기본적으로",
  ")"
)

$full = [string]::Join("", $parts)
$p.Range.InsertAfter($full)

# Compute the absolute document positions of each segment boundary so we
# can force Word to keep the segments as separate runs (runs with
# identical formatting otherwise get coalesced back into a single run).
$offset = $start0
$boundaries = @()
foreach ($part in $parts) {
    $offset = $offset + $part.Length
    $boundaries += $offset
}
$splitPoints = $boundaries[0..($boundaries.Length - 2)]
$paraEnd = $p.Range.End - 1

# Toggling formatting on the tail of the paragraph (from the boundary to
# the end of the text) and immediately reverting it forces a run split at
# that boundary without altering the final (empty) run formatting.
# Walk right-to-left so earlier splits don't shift later boundaries.
for ($i = $splitPoints.Length - 1; $i -ge 0; $i--) {
    $pos = $splitPoints[$i]
    $rr = $d.Range($pos, $paraEnd)
    $rr.Font.Bold = 1
    $rr.Font.Bold = 0
}

# Append a new blank paragraph, then a paragraph containing a single
# space, matching the trailing structure added at the end of the body.
$blankPara = $d.Paragraphs.Add()
$spacePara = $d.Paragraphs.Add()
$spacePara.Range.InsertAfter(" ")
